# Change the Date column (A2:A6) from text dates to real Excel date
# values formatted as m/dd/yyyy (commit: "change excel dates format to
# m/dd/yyyy"). Also append a new (currently empty) dated row 7 below the
# existing data, matching the author's next step of continuing the
# calendar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    (Get-Date -Year 2026 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0),
    (Get-Date -Year 2026 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0),
    (Get-Date -Year 2026 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0),
    (Get-Date -Year 2026 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0),
    (Get-Date -Year 2026 -Month 2 -Day 27 -Hour 0 -Minute 0 -Second 0)
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "m/dd/yyyy"
}

# New row 7: the next day in the sequence, date column only (no other
# data filled in yet), carrying the same date style.
$newDateCell = $ws.Cells.Item(7, 1)
$newDateCell.NumberFormat = "m/dd/yyyy"

# Move the active selection down to the newly added row, ready for entry.
$ws.Range("A7").Select()

# Refresh the ignoredErrors ("numbers stored as text") ranges now that
# column A holds real numeric dates instead of text - the warning no
# longer applies to column A, but the other previously-unflagged header/
# data cells keep the same suppression as before.
$ws.Range("A1:G1,B6:G6,G4,G5,B3:G3,B2:G2").ErrorCheckingOptions.NumberAsText = $false
